$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph dynamically
# (avoid piping COM objects through Where-Object/Select-Object, which
# breaks live property access on this host; use a plain indexed loop).
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "Docente(s)*") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Docente(s) Responsável(eis)' paragraph"
}

$headingPara = $d.Paragraphs.Item($targetIndex)
$headingRange = $headingPara.Range
$anchor = $d.Range($headingRange.Start, $headingRange.End)

# Create a fresh empty paragraph right after the heading ...
$anchor.InsertParagraphAfter() | Out-Null

# ... and fill it with the two-run list item (name + line break, then the
# second name) by injecting the exact WordprocessingML for the paragraph.
# This keeps the two names as separate <w:r> runs, matching how Word
# represents text that has a manual line break between differently-typed
# spans, rather than collapsing them into a single run.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range

$xml = '<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>1176388 - Luiz Tadeu Fernandes Eleno</w:t><w:br/></w:r><w:r><w:t>7797767 - Viktor Pastoukhov</w:t></w:r></w:p>'
$newRange.InsertXML($xml) | Out-Null
